$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.69
$ws.Cells.Item(2, 8).Value = 0.69
$ws.Cells.Item(3, 7).Value = 0.0033
$ws.Cells.Item(3, 8).Value = 0.0066
$ws.Cells.Item(4, 7).Value = 0.0077
$ws.Cells.Item(4, 8).Value = 0.0693
$ws.Cells.Item(5, 7).Value = 0.0737
$ws.Cells.Item(5, 8).Value = 0.0737
$ws.Cells.Item(6, 7).Value = 0.0052
$ws.Cells.Item(6, 8).Value = 0.026
$ws.Cells.Item(7, 7).Value = 0.1036
$ws.Cells.Item(7, 8).Value = 0.1036
$ws.Cells.Item(8, 7).Value = 0.0983
$ws.Cells.Item(8, 8).Value = 0.0983
$ws.Cells.Item(10, 7).Value = 0.0144
$ws.Cells.Item(10, 8).Value = 0.0144
$ws.Cells.Item(11, 7).Value = 0.0107
$ws.Cells.Item(11, 8).Value = 0.0107
$ws.Cells.Item(12, 7).Value = 49.32
$ws.Cells.Item(12, 8).Value = 49.32
$ws.Cells.Item(13, 1).Value = 'CDBA5150-HF'
$ws.Cells.Item(13, 2).Value = 'No Description Available'
$ws.Cells.Item(14, 1).Value = 'B2B-EH-A (LF)(SN)'
$ws.Cells.Item(14, 3).Value = 'B2B-EH-A(LF)(SN)'
$ws.Cells.Item(14, 4).Value = 5.0
$ws.Cells.Item(14, 5).Value = 'CONN_B2B-EH-A (LF)(SN)_JST'
$ws.Cells.Item(14, 7).Value = 0.14
$ws.Cells.Item(14, 8).Value = 0.7
$ws.Cells.Item(15, 1).Value = 'B3B-EH-A(LF)(SN)'
$ws.Cells.Item(15, 3).Value = 'B3B-EH-A(LF)(SN)'
$ws.Cells.Item(15, 4).Value = 3.0
$ws.Cells.Item(15, 5).Value = 'CONN3_B3B-EH-A(LF)(SN)_JST'
$ws.Cells.Item(15, 7).Value = 0.19
$ws.Cells.Item(15, 8).Value = 0.57
$ws.Cells.Item(16, 1).Value = 'B4B-XH-A (LF)(SN)'
$ws.Cells.Item(16, 3).Value = 'B4B-XH-A(LF)(SN)'
$ws.Cells.Item(16, 4).Value = 1.0
$ws.Cells.Item(16, 5).Value = 'CONN_B4B-XH-A (LF)(SN)_JST'
$ws.Cells.Item(16, 7).Value = 0.21
$ws.Cells.Item(16, 8).Value = 0.21
$ws.Cells.Item(17, 1).Value = 74438323100.0
$ws.Cells.Item(17, 2).Value = 'Power Multilayer Inductor WE-MAPI, L=10.0µH'
$ws.Cells.Item(17, 3).Value = 74438323100.0
$ws.Cells.Item(17, 5).Value = 'SMD 2510'
$ws.Cells.Item(17, 7).Value = 1.13
$ws.Cells.Item(17, 8).Value = 1.13
$ws.Cells.Item(19, 4).Value = 3.0
$ws.Cells.Item(20, 1).Value = 'Header 5'
$ws.Cells.Item(20, 2).Value = 'Header, 5-Pin'
$ws.Cells.Item(20, 3).Value = 'B5B-XH-A-(LF)(SN)'
$ws.Cells.Item(20, 4).Value = 2.0
$ws.Cells.Item(20, 5).Value = 'CONN_B5B-XH-A (LF)(SN)_JST'
$ws.Cells.Item(20, 7).Value = 0.268
$ws.Cells.Item(20, 8).Value = 0.536
$ws.Cells.Item(21, 3).Value = 'CRCW0603100MJPEAHR'
$ws.Cells.Item(21, 4).Value = 1.0
$ws.Cells.Item(21, 5).Value = '1608[0603]'
$ws.Cells.Item(21, 6).Value = '100M'
$ws.Cells.Item(21, 7).Value = 0.2552
$ws.Cells.Item(21, 8).Value = 0.2552
$ws.Cells.Item(22, 3).Value = 'RC0402FR-071KL'
$ws.Cells.Item(22, 4).Value = 4.0
$ws.Cells.Item(22, 5).Value = '0402-A'
$ws.Cells.Item(22, 6).Value = '1k'
$ws.Cells.Item(22, 7).Value = 0.012
$ws.Cells.Item(22, 8).Value = 0.12
$ws.Cells.Item(23, 3).Value = 'ESR10EZPJ681'
$ws.Cells.Item(23, 5).Value = '6-0805_M'
$ws.Cells.Item(23, 6).Value = 680.0
$ws.Cells.Item(23, 7).Value = 0.1
$ws.Cells.Item(23, 8).Value = 0.1
$ws.Cells.Item(24, 3).Value = 'CR0603-JW-331ELF'
$ws.Cells.Item(24, 4).Value = 1.0
$ws.Cells.Item(24, 5).Value = '1608[0603]'
$ws.Cells.Item(24, 6).Value = 330.0
$ws.Cells.Item(24, 7).Value = 0.003
$ws.Cells.Item(24, 8).Value = 0.003
$ws.Cells.Item(25, 3).Value = 'RC0402FR-0710KL'
$ws.Cells.Item(25, 4).Value = 2.0
$ws.Cells.Item(25, 6).Value = '10k'
$ws.Cells.Item(25, 7).Value = 0.059
$ws.Cells.Item(25, 8).Value = 0.118
$ws.Cells.Item(26, 3).Value = 'MCS04020C1004FE000'
$ws.Cells.Item(26, 4).Value = 4.0
$ws.Cells.Item(26, 5).Value = '0402-A'
$ws.Cells.Item(26, 6).Value = '1M'
$ws.Cells.Item(26, 7).Value = 0.035
$ws.Cells.Item(26, 8).Value = 0.14
$ws.Cells.Item(27, 3).Value = 'RR0816P-513-D'
$ws.Cells.Item(27, 4).Value = 2.0
$ws.Cells.Item(27, 5).Value = '1608[0603]'
$ws.Cells.Item(27, 6).Value = '51k'
$ws.Cells.Item(27, 7).Value = 0.1
$ws.Cells.Item(27, 8).Value = 0.2
$ws.Cells.Item(28, 1).Value = 'Res1'
$ws.Cells.Item(28, 2).Value = 'Resistor'
$ws.Cells.Item(28, 3).Value = 'RC0402JR-070RL'
$ws.Cells.Item(28, 5).Value = '0402-A'
$ws.Cells.Item(28, 6).Value = 0.0
$ws.Cells.Item(28, 7).Value = 0.0016
$ws.Cells.Item(28, 8).Value = 0.0016
$ws.Cells.Item(29, 1).Value = 'LMR62014XMF/NOPD'
$ws.Cells.Item(29, 3).Value = 'LMR62014XMF/NOPB'
$ws.Cells.Item(29, 5).Value = 'SOT95P280X145-5N'
$ws.Cells.Item(29, 7).Value = 0.81765
$ws.Cells.Item(29, 8).Value = 0.81765
$ws.Cells.Item(30, 1).Value = 'SSCMNNN030PA2A3'
$ws.Cells.Item(30, 2).Value = 'No Description Available'
$ws.Cells.Item(30, 3).Value = 'SSCMNNN030PA2A3'
$ws.Cells.Item(30, 4).Value = 2.0
$ws.Cells.Item(30, 5).Value = 'SSCMNN_HNW-L'
$ws.Cells.Item(30, 7).Value = 35.0
$ws.Cells.Item(30, 8).Value = 70.0
$ws.Cells.Item(31, 1).Value = 'PIC18LF45K50-I/PT'
$ws.Cells.Item(31, 2).Value = 'Low-Power, High-Performance Microcontroller with XLP Technology, 48 MHz, 36 I/O, -40 to 85 degC, 44-pin TQFP (PT44), Tray'
$ws.Cells.Item(31, 3).Value = 'PIC18LF45K50-I/PT'
$ws.Cells.Item(31, 5).Value = 'MCHP-TQFP-PT44_L'
$ws.Cells.Item(31, 7).Value = 2.45
$ws.Cells.Item(31, 8).Value = 2.45
$ws.Cells.Item(32, 1).Value = 'OPA336NA/250'
$ws.Cells.Item(32, 2).Value = 'Single-Supply, MicroPower CMOS Operational Amplifier MicroAmplifier(TM) Series, 2.3 to 5.5 V, -55 to 125 degC, 5-pin SOT23 (DBV5), Green (RoHS & no Sb/Br)'
$ws.Cells.Item(32, 3).Value = 'OPA336NA/250'
$ws.Cells.Item(32, 4).Value = 9.0
$ws.Cells.Item(32, 5).Value = 'DBV0005A_M'
$ws.Cells.Item(32, 7).Value = 1.57
$ws.Cells.Item(32, 8).Value = 15.7
$ws.Cells.Item(33, 1).Value = 'MCP2515-I/SO'
$ws.Cells.Item(33, 2).Value = 'Stand-Alone CAN Controller With SPI Interface, 18-Pin SOIC, Industrial Temperature'
$ws.Cells.Item(33, 3).Value = 'MCP2515-I/SO'
$ws.Cells.Item(33, 5).Value = 'SOIC-SO18_L'
$ws.Cells.Item(33, 7).Value = 1.79
$ws.Cells.Item(33, 8).Value = 1.79
$ws.Cells.Item(34, 1).Value = 'TCAN334DR'
$ws.Cells.Item(34, 2).Value = 'No Description Available'
$ws.Cells.Item(34, 3).Value = 'TCAN334DR'
$ws.Cells.Item(34, 5).Value = 'D0008A_N'
$ws.Cells.Item(34, 6).Value = 'Value'
$ws.Cells.Item(34, 7).Value = 2.3
$ws.Cells.Item(34, 8).Value = 2.3
$ws.Cells.Item(35, 1).Value = 'TMP235A2DBZR'
$ws.Cells.Item(35, 2).Value = 'LOW-POWER HIGH-ACCURACY ANALOG O'
$ws.Cells.Item(35, 3).Value = 'TMP235A2DBZR'
$ws.Cells.Item(35, 5).Value = 'FP-DBZ0003A-MFG'
$ws.Cells.Item(35, 7).Value = 0.7584
$ws.Cells.Item(35, 8).Value = 0.7584
$ws.Cells.Item(36, 1).Value = 'D5V0L2B3W-7'
$ws.Cells.Item(36, 2).Value = 'No Description Available'
$ws.Cells.Item(36, 3).Value = 'D5V0L2B3W-7'
$ws.Cells.Item(36, 5).Value = 'SOT-323_DIO'
$ws.Cells.Item(36, 7).Value = 0.2301
$ws.Cells.Item(36, 8).Value = 0.2301
$ws.Cells.Item(37, 1).Value = 'LSM6DSMTR'
$ws.Cells.Item(37, 2).Value = 'IMU ACCEL/GYRO I2C/SPI 14LGA'
$ws.Cells.Item(37, 3).Value = 'LSM6DSMTR'
$ws.Cells.Item(37, 4).Value = 1.0
$ws.Cells.Item(37, 5).Value = 'FP-LGA-14L-DM00249496-MFG'
$ws.Cells.Item(37, 7).Value = 2.78
$ws.Cells.Item(37, 8).Value = 2.78
$ws.Cells.Item(20, 6).ClearContents()
$ws.Cells.Item(29, 2).ClearContents()
$ws.Cells.Item(33, 6).ClearContents()
